$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.885.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.321.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.68%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.319.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "659.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.864.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.818.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.321.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -2.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("E23").Value = "  +4.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.89"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "567.49"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.89%  "
$ws.Range("E33").Value = "  -0.76%  "
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.676.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.89%  "
$ws.Range("E40").Value = "  +0.41%  "
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("E42").Value = "  -5.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  -2.74%  "
$ws.Range("E51").Value = "  +9.52%  "
